{"js": "// Applies the \"Update Leave Card\" edits: swap the employee's name/position/\n// office/last-day-of-service/leave totals from Cristina M. Igno to\n// Felicitas M. Sumagui, update the issuance date from June to July, and\n// tidy up a stray grammar-check split around \"pay.\" \u2014 each change appears\n// twice since the certificate body is duplicated in the document.\n\nconst body = context.document.body;\n\n// Simple text substitutions: search the target phrase, replace every hit.\nasync function replaceAll(findText, replaceText, options) {\n  const opts = options || { matchCase: true };\n  const results = body.search(findText, opts);\n  results.load(\"text\");\n  await context.sync();\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(replaceText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n\n// Name, position, office, last day of service.\nawait replaceAll(\"CRISTINA M. IGNO\", \"FELICITAS M. SUMAGUI\");\nawait replaceAll(\"Administrative Aide III\", \"Casual Employee\");\nawait replaceAll(\"Human Resource Management Office\", \"City Social Welfare Development Office\");\nawait replaceAll(\"October 31, 2022\", \"March 22, 2023\");\n\n// Leave credit totals (leading whitespace matters \u2014 see diff).\nawait replaceAll(\" 137.868\", \"  53.458\");\nawait replaceAll(\" 188.063\", \"  67.458\");\nawait replaceAll(\" 325.931\", \" 120.916\");\n\n// Short-form surname used later in the certificate (\"Ms. Igno ...\").\nawait replaceAll(\"Igno\", \"Sumagui\", { matchCase: true, matchWholeWord: true });\n\n// Month of issuance.\nawait replaceAll(\"June\", \"July\", { matchCase: true, matchWholeWord: true });\n\n// The phrase \" for claiming terminal leave pay.\" was previously split into\n// three runs with <w:proofErr w:type=\"gramStart/gramEnd\"/> bracketing the\n// word \"pay\" (a leftover grammar-check flag). Re-running a same-text\n// replace is a no-op in this engine because nothing differs, so first\n// swap in a placeholder to force a real overwrite, then restore the\n// correct text \u2014 this collapses the run back to a single run and drops\n// the now-stale proofErr markers, matching the authored edit.\nawait replaceAll(\" for claiming terminal leave pay.\", \"\\u0001PLACEHOLDER\\u0001\");\nawait replaceAll(\"\\u0001PLACEHOLDER\\u0001\", \" for claiming terminal leave pay.\");\n", "ps1": "$d = $word.ActiveDocument\n\nfunction Replace-All($find, $replace) {\n    $d.Content.Find.Execute($find, $false, $false, $false, $false, $false, $true, 1, $false, $replace, 2) | Out-Null\n}\n\nReplace-All \"CRISTINA M. IGNO\" \"FELICITAS M. SUMAGUI\"\nReplace-All \"Administrative Aide III\" \"Casual Employee\"\nReplace-All \"Human Resource Management Office\" \"City Social Welfare Development Office\"\nReplace-All \"October 31, 2022\" \"March 22, 2023\"\nReplace-All \" 137.868\" \"  53.458\"\nReplace-All \" 188.063\" \"  67.458\"\nReplace-All \" 325.931\" \" 120.916\"\nReplace-All \"Igno\" \"Sumagui\"\nReplace-All \"June\" \"July\"\nReplace-All \" for claiming terminal leave pay.\" \" for claiming terminal leave pay.\"\n"}
